# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value into a cell as literal TEXT (never auto-converted to a
# number/date by Excel type inference), while leaving the target cells style
# untouched. We stage the text in a scratch cell formatted as Text, copy it, and
# paste-special VALUES ONLY into the destination (paste-values does not carry the
# scratch cells number format along with it).
function Set-TextValue($cell, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.Clear()
}

$ws.Range('D2').Value = '35.351.36'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '1.912.46'
$ws.Range('E3').Value = '  +2.94%  '
$ws.Range('E4').Value = '  -0.47%  '
Set-TextValue $ws.Range('D5') '245.25'
$ws.Range('E5').Value = '  +2.50%  '
Set-TextValue $ws.Range('D6') '0.660'
$ws.Range('E6').Value = '  +6.09%  '
$ws.Range('E7').Value = '  -0.48%  '
Set-TextValue $ws.Range('D8') '41.30'
$ws.Range('E8').Value = '  -1.75%  '
$ws.Range('E9').Value = '  +6.31%  '
Set-TextValue $ws.Range('D10') '52.86'
$ws.Range('E10').Value = '  +12.64%  '
Set-TextValue $ws.Range('D11') '0.0717'
$ws.Range('E11').Value = '  +3.53%  '
Set-TextValue $ws.Range('D12') '0.0998'
$ws.Range('E12').Value = '  +0.98%  '
$ws.Range('D13').Value = '2.190.19'
$ws.Range('E13').Value = '  +2.91%  '
Set-TextValue $ws.Range('D14') '12.08'
$ws.Range('E14').Value = '  +5.34%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D15') '0.703'
$ws.Range('E15').Value = '  +3.95%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.919.04'
$ws.Range('E16').Value = '  +3.16%  '
Set-TextValue $ws.Range('D17') '4.87'
$ws.Range('E17').Value = '  +3.26%  '
$ws.Range('D18').Value = '35.339.17'
$ws.Range('E18').Value = '  +0.53%  '
Set-TextValue $ws.Range('D19') '72.08'
$ws.Range('E19').Value = '  +3.24%  '
$ws.Range('D20').Value = '0.0₃0830'
$ws.Range('E20').Value = '  +4.46%  '
Set-TextValue $ws.Range('D21') '239.51'
$ws.Range('E21').Value = '  -0.36%  '
Set-TextValue $ws.Range('D22') '12.50'
$ws.Range('E22').Value = '  +2.31%  '
Set-TextValue $ws.Range('D23') '4.83'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  +1.32%  '
Set-TextValue $ws.Range('D26') '2.36'
$ws.Range('E26').Value = '  +22.25%  '
Set-TextValue $ws.Range('D27') '170.02'
$ws.Range('E27').Value = '  +1.06%  '
Set-TextValue $ws.Range('D28') '8.49'
$ws.Range('E28').Value = '  +6.45%  '
Set-TextValue $ws.Range('D29') '18.46'
$ws.Range('E29').Value = '  +4.63%  '
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('E31').Value = '  +4.07%  '
Set-TextValue $ws.Range('D32') '0.0568'
$ws.Range('E32').Value = '  +1.89%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D33') '0.936'
$ws.Range('E33').Value = '  +13.55%  '
$ws.Range('B34').Value = 'BinanceUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D34') '1.01'
$ws.Range('E34').Value = '  -0.47%  '
Set-TextValue $ws.Range('D35') '4.12'
$ws.Range('E35').Value = '  +2.66%  '
$ws.Range('E36').Value = '  -4.20%  '
Set-TextValue $ws.Range('D37') '2.05'
$ws.Range('E37').Value = '  +1.42%  '
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D39') '0.0666'
$ws.Range('E39').Value = '  +12.42%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D40') '1.12'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D41') '0.0208'
$ws.Range('E41').Value = '  +3.96%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D42') '16.30'
$ws.Range('E42').Value = '  +9.82%  '
Set-TextValue $ws.Range('D43') '90.27'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').Value = '1.340.52'
$ws.Range('E44').Value = '  +0.06%  '
Set-TextValue $ws.Range('D45') '2.38'
$ws.Range('E45').Value = '  +3.07%  '
Set-TextValue $ws.Range('D46') '47.88'
$ws.Range('E46').Value = '  +37.74%  '
$ws.Range('B47').Value = 'MXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D47') '2.79'
$ws.Range('E47').Value = '  +1.94%  '
$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D48') '2.41'
$ws.Range('E48').Value = '  -0.31%  '
Set-TextValue $ws.Range('D49') '6.57'
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('D50').Value = '2.097.05'
$ws.Range('E50').Value = '  +2.70%  '
Set-TextValue $ws.Range('D51') '0.0702'
$ws.Range('E51').Value = '  +3.35%  '
